$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: the paragraph that contains only the "_GoBack" bookmark becomes the
# opening Jinja tag:
#   {% if special_bond_conditions.admin_license_suspension_ordered is true %}
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$insertPos = $bm.Range.Start
$r1 = $d.Range($insertPos, $insertPos)

$part1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>special_bond_conditions.admin_license_suspension_ordered</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> is true %}</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r1.InsertXML($part1)

# The old _GoBack bookmark is still sitting right after the text we just
# inserted (same paragraph) - remove it here; it gets re-created below, right
# after the Administrative License Suspension paragraph's last sentence.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# Part 2: append the closing tag "{% endif %}" right after "...IS NOT STAYED
# ." and re-create the "_GoBack" bookmark immediately after it.
# ---------------------------------------------------------------------------
$findRange = $d.Content
[void]$findRange.Find.Execute("IS NOT STAYED .", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Grab the whole paragraph so we can replace it in one shot - replacing a
# range that spans only part of a paragraph's tail (but not all of it) can
# split the paragraph instead of appending in place, so we rebuild the full
# paragraph content instead.
$targetPara = $findRange.Paragraphs(1)
$paraXml = $targetPara.Range.WordOpenXML
$null = $paraXml -match '(?s)<w:body>.*?<w:p\b[^>]*>(.*?)</w:p>'
$originalInner = $matches[1]

$newRuns = '<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$fullInner = $originalInner + $newRuns

$part2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p>' + $fullInner + '</w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2 = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
[void]$r2.InsertXML($part2)

Write-Host "edit complete"
